$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "r775"
$ws.Range("B20").Value = "cameron"
$ws.Range("C20").Value = "that feeling when you `"had it`" 5 tries ago"
$ws.Range("D20").Value = "2025-10-01 16:20:33"
